$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5049.8335
$ws.Range("I64").Value = 5374.75
$ws.Range("K64").Value = 5374.75
$ws.Range("M64").Value = -5126.75
$ws.Range("H67").Value = 5049.8335
$ws.Range("I67").Value = 5374.75
$ws.Range("K67").Value = 5374.75
$ws.Range("M67").Value = -4516.75
$ws.Range("H70").Value = 1432.6
$ws.Range("I70").Value = 1057.1428
$ws.Range("J70").Value = 1761.125
$ws.Range("K70").Value = 3171.4284
$ws.Range("L70").Value = 5283.375
$ws.Range("M70").Value = -2901.4284
$ws.Range("N70").Value = -5823.375
$ws.Range("H73").Value = 1432.6
$ws.Range("I73").Value = 1057.1428
$ws.Range("J73").Value = 1761.125
$ws.Range("K73").Value = 3171.4284
$ws.Range("L73").Value = 5283.375
$ws.Range("M73").Value = -2235.4284
$ws.Range("N73").Value = -7155.375
$ws.Range("H106").Value = 1855.875
$ws.Range("I106").Value = 1692.4286
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1692.4286
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1061.4286
$ws.Range("N106").Value = -4262
$ws.Range("H141").Value = 809.1429000000001
$ws.Range("I141").Value = 799.6
$ws.Range("K141").Value = 2398.8
$ws.Range("M141").Value = 2781.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1800.6945
$ws.Range("I2").Value = 1870.1538
$ws.Range("K2").Value = 1870.1538
$ws.Range("M2").Value = -1757.1538
$ws.Range("H32").Value = 2377
$ws.Range("I32").Value = 2377
$ws.Range("K32").Value = 2377
$ws.Range("M32").Value = -2090
$ws.Range("H45").Value = 2028.5714
$ws.Range("I45").Value = 2400
$ws.Range("J45").Value = 1657.1428
$ws.Range("K45").Value = 2400
$ws.Range("L45").Value = 1657.1428
$ws.Range("M45").Value = -2023
$ws.Range("N45").Value = -2411.1428
$ws.Range("H97").Value = 2421
$ws.Range("I97").Value = 3076.65
$ws.Range("K97").Value = 3076.65
$ws.Range("M97").Value = -2580.65
$ws.Range("H116").Value = 1800.6945
$ws.Range("I116").Value = 1870.1538
$ws.Range("K116").Value = 1870.1538
$ws.Range("M116").Value = 423.8462
$ws.Range("H132").Value = 1648.6
$ws.Range("I132").Value = 1276
$ws.Range("J132").Value = 2518
$ws.Range("K132").Value = 3828
$ws.Range("L132").Value = 7554
$ws.Range("M132").Value = -1298
$ws.Range("N132").Value = -12614
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1800.6945
$ws.Range("I3").Value = 1870.1538
$ws.Range("K3").Value = 1870.1538
$ws.Range("M3").Value = -1756.1538
$ws.Range("H94").Value = 495.73914
$ws.Range("I94").Value = 381.9091
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 381.9091
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = 69.09089999999998
$ws.Range("N94").Value = -3902
$ws.Range("H132").Value = 700000
$ws.Range("J132").Value = 700000
$ws.Range("L132").Value = 700000
$ws.Range("N132").Value = -710120
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125
$ws.Range("I7").Value = 125
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 125
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -12
$ws.Range("N7").ClearContents()
$ws.Range("H62").Value = 5510.8237
$ws.Range("I62").Value = 5553.1816
$ws.Range("J62").Value = 5433.1665
$ws.Range("K62").Value = 5553.1816
$ws.Range("L62").Value = 5433.1665
$ws.Range("M62").Value = -4929.1816
$ws.Range("N62").Value = -6681.1665
$ws.Range("H65").Value = 5510.8237
$ws.Range("I65").Value = 5553.1816
$ws.Range("J65").Value = 5433.1665
$ws.Range("K65").Value = 27765.908
$ws.Range("L65").Value = 27165.8325
$ws.Range("M65").Value = -24645.908
$ws.Range("N65").Value = -33405.8325
$ws.Range("H69").Value = 28000
$ws.Range("I69").Value = 28000
$ws.Range("K69").Value = 28000
$ws.Range("M69").Value = -27251
$ws.Range("H72").Value = 28000
$ws.Range("I72").Value = 28000
$ws.Range("K72").Value = 84000
$ws.Range("M72").Value = -80256
$ws.Range("H86").Value = 2977.8333
$ws.Range("I86").Value = 1999.5
$ws.Range("J86").Value = 3467
$ws.Range("K86").Value = 1999.5
$ws.Range("L86").Value = 3467
$ws.Range("M86").Value = -876.5
$ws.Range("N86").Value = -5713
$ws.Range("H89").Value = 2977.8333
$ws.Range("I89").Value = 1999.5
$ws.Range("J89").Value = 3467
$ws.Range("K89").Value = 9997.5
$ws.Range("L89").Value = 17335
$ws.Range("M89").Value = -4381.5
$ws.Range("N89").Value = -28567
$ws.Range("H132").Value = 1059.4524
$ws.Range("I132").Value = 859.4
$ws.Range("J132").Value = 2059.7144
$ws.Range("K132").Value = 2578.2
$ws.Range("L132").Value = 6179.1432
$ws.Range("M132").Value = -48.19999999999982
$ws.Range("N132").Value = -11239.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3921749
$ws.Range("I2").Value = 5882592.5
$ws.Range("J2").Value = 61.4
$ws.Range("K2").Value = 35295555
$ws.Range("L2").Value = 368.4
$ws.Range("M2").Value = -35295442
$ws.Range("N2").Value = -594.4
$ws.Range("H137").Value = 70818.92999999999
$ws.Range("J137").Value = 341833.34
$ws.Range("L137").Value = 1025500.02
$ws.Range("N137").Value = -1035700.02
$ws.Range("H138").Value = 1923.8889
$ws.Range("I138").Value = 1452.8572
$ws.Range("J138").Value = 3572.5
$ws.Range("K138").Value = 4358.571599999999
$ws.Range("L138").Value = 10717.5
$ws.Range("M138").Value = 781.4284000000007
$ws.Range("N138").Value = -20997.5
$ws.Range("H139").Value = 34729.035
$ws.Range("I139").Value = 39553.5
$ws.Range("J139").Value = 3370
$ws.Range("K139").Value = 118660.5
$ws.Range("L139").Value = 10110
$ws.Range("M139").Value = -113520.5
$ws.Range("N139").Value = -20390
$ws.Range("H140").Value = 217879.28
$ws.Range("I140").Value = 276446.38
$ws.Range("J140").Value = 3133.3333
$ws.Range("K140").Value = 829339.14
$ws.Range("L140").Value = 9399.999899999999
$ws.Range("M140").Value = -824159.14
$ws.Range("N140").Value = -19759.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 31428.572
$ws.Range("J53").Value = 31428.572
$ws.Range("L53").Value = 31428.572
$ws.Range("N53").Value = -32464.572
$ws.Range("H93").Value = 1798.2727
$ws.Range("I93").Value = 1281.2
$ws.Range("J93").Value = 2478.6316
$ws.Range("K93").Value = 1281.2
$ws.Range("L93").Value = 2478.6316
$ws.Range("M93").Value = -33.20000000000005
$ws.Range("N93").Value = -4974.631600000001
$ws.Range("H100").Value = 3832931.2
$ws.Range("I100").Value = 4631017
$ws.Range("J100").Value = 2120
$ws.Range("K100").Value = 4631017
$ws.Range("L100").Value = 2120
$ws.Range("M100").Value = -4630476
$ws.Range("N100").Value = -3202
$ws.Range("H136").Value = 3381.6323
$ws.Range("I136").Value = 1833.6086
$ws.Range("K136").Value = 5500.825800000001
$ws.Range("M136").Value = -2950.825800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 839.4
$ws.Range("I126").Value = 843.7778
$ws.Range("K126").Value = 2531.3334
$ws.Range("M126").Value = -61.33339999999998
$ws.Range("H132").Value = 1023.22534
$ws.Range("I132").Value = 782.2549
$ws.Range("J132").Value = 1637.7
$ws.Range("K132").Value = 2346.7647
$ws.Range("L132").Value = 4913.1
$ws.Range("M132").Value = 183.2352999999998
$ws.Range("N132").Value = -9973.1
